$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.563.39'
$ws.Range('E2').Value = '  +5.79%  '
$ws.Range('D3').Value = '3.253.83'
$ws.Range('E3').Value = '  +3.81%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.00'
$ws.Range('E5').Value = '  +2.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.27'
$ws.Range('E6').Value = '  +7.85%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  -2.35%  '
$ws.Range('D9').Value = '3.255.20'
$ws.Range('E9').Value = '  +4.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.130'
$ws.Range('E10').Value = '  +5.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.79'
$ws.Range('E11').Value = '  +3.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.412'
$ws.Range('E12').Value = '  +5.49%  '
$ws.Range('D13').Value = '3.801.67'
$ws.Range('E13').Value = '  +3.41%  '
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.93'
$ws.Range('E15').Value = '  +3.64%  '
$ws.Range('D16').Value = '67.617.24'
$ws.Range('E16').Value = '  +6.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000168'
$ws.Range('E17').Value = '  +3.81%  '
$ws.Range('D18').Value = '3.257.75'
$ws.Range('E18').Value = '  +4.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.82'
$ws.Range('E19').Value = '  +2.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.42'
$ws.Range('E20').Value = '  +4.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '373.55'
$ws.Range('E21').Value = '  +6.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.59'
$ws.Range('E22').Value = '  +6.31%  '
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.07'
$ws.Range('E24').Value = '  +4.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.510'
$ws.Range('E25').Value = '  +2.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000119'
$ws.Range('E26').Value = '  +4.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.59'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('E28').Value = '  +3.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.97'
$ws.Range('E30').Value = '  +5.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.65'
$ws.Range('E31').Value = '  +4.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.58'
$ws.Range('E32').Value = '  +3.60%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.27'
$ws.Range('E34').Value = '  +6.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.87'
$ws.Range('E35').Value = '  +4.60%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.50'
$ws.Range('E36').Value = '  +5.52%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.34'
$ws.Range('E37').Value = '  +5.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.859'
$ws.Range('E38').Value = '  +5.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.85'
$ws.Range('E39').Value = '  +10.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.81'
$ws.Range('E40').Value = '  +14.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.80'
$ws.Range('E41').Value = '  +2.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.59'
$ws.Range('E42').Value = '  +7.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '361.57'
$ws.Range('E43').Value = '  +15.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.39'
$ws.Range('E44').Value = '  +6.27%  '
$ws.Range('D45').Value = '2.707.45'
$ws.Range('E45').Value = '  +2.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.63'
$ws.Range('E46').Value = '  +8.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.59'
$ws.Range('E47').Value = '  +3.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0674'
$ws.Range('E48').Value = '  +4.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0278'
$ws.Range('E49').Value = '  +3.70%  '
$ws.Range('E50').Value = '  +8.45%  '
$ws.Range('E51').Value = '  +1.93%  '
